$d = $word.ActiveDocument

$oldUrl = "https://github.com/EduriMaryJones/Pollen-sProfiling_Project.git"
$newUrl = "https://github.com/EduriMaryJones/Pollen-s-Profiling-Automated-Classification-of-Pollen-Grains.git"

# Locate the run that holds the old repository URL and rewrite its text.
$text = $d.Content.Text
$startUrl = $text.IndexOf($oldUrl)
if ($startUrl -lt 0) {
    throw "Could not locate the old repository URL in the document."
}
$endUrl = $startUrl + $oldUrl.Length

$urlRange = $d.Range($startUrl, $endUrl)
$urlRange.Text = $newUrl

# The URL run used to be followed by a separate, Arial-formatted run that
# contained a single trailing space. That run must disappear entirely.
$text = $d.Content.Text
$afterUrl = $text.IndexOf($newUrl) + $newUrl.Length
if ($afterUrl -lt $text.Length -and $text.Substring($afterUrl, 1) -eq " ") {
    $trailingSpace = $d.Range($afterUrl, $afterUrl + 1)
    $trailingSpace.Delete()
}
